# SE-2873 Extended to use notional cost
# Adds a "notional_amount" column (M) to the transactions sheet, converts the
# futures-style transaction in row 5 to an "OpenContract" txn_type whose
# total_consideration becomes the contract count (50) while the previous
# total_consideration value (746999) moves to the new notional_amount column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("transactions")

# Row 5 (txnid_0004): txn_type -> OpenContract, total_consideration -> 50,
# notional_amount -> 746999 (the old total_consideration value).
# "OpenContract" is written first so it lands before "notional_amount" in
# the shared-strings table, matching the saved file.
$ws.Range("C5").Value = "OpenContract"

# New header for column M
$ws.Range("M1").Value = "notional_amount"
$ws.Range("M1").HorizontalAlignment = -4108

$ws.Range("L5").Value = 50
$ws.Range("M5").Value = 746999

# Populate notional_amount with 0 for all other transaction rows (2-22)
$ws.Range("M2").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("M8").Value = 0
$ws.Range("M9").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("M15").Value = 0
$ws.Range("M16").Value = 0
$ws.Range("M17").Value = 0
$ws.Range("M18").Value = 0
$ws.Range("M19").Value = 0
$ws.Range("M20").Value = 0
$ws.Range("M21").Value = 0
$ws.Range("M22").Value = 0

# Re-activate the transactions sheet and restore the selection to the cell
# just past the new last column/row of data, matching the saved view state.
$ws.Activate() | Out-Null
$ws.Range("M23").Select() | Out-Null
